$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its text formatting so values such as
# "313.17" or "1.000" are not reinterpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.650.90"
$ws.Range("E2").Value = "  -4.27%  "
$ws.Range("D3").Value = "1.846.33"
$ws.Range("E3").Value = "  -3.66%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "313.17"
$ws.Range("E5").Value = "  -3.46%  "
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").Value = "0.4252"
$ws.Range("E7").Value = "  -6.79%  "
$ws.Range("D8").Value = "0.3639"
$ws.Range("E8").Value = "  -4.16%  "
$ws.Range("D9").Value = "43.77"
$ws.Range("D10").Value = "0.07219"
$ws.Range("E10").Value = "  -6.69%  "
$ws.Range("D11").Value = "0.8991"
$ws.Range("E11").Value = "  -7.73%  "
$ws.Range("D12").Value = "20.68"
$ws.Range("E12").Value = "  -7.06%  "
$ws.Range("D13").Value = "1.815.44"
$ws.Range("E13").Value = "  -5.02%  "
$ws.Range("D14").Value = "6.581"
$ws.Range("E14").Value = "  -5.52%  "
$ws.Range("D15").Value = "5.351"
$ws.Range("E15").Value = "  -5.93%  "
$ws.Range("D16").Value = "0.06802"
$ws.Range("E16").Value = "  -2.66%  "
$ws.Range("D18").Value = "77.56"
$ws.Range("E18").Value = "  -8.01%  "
$ws.Range("D19").Value = "0.000008875"
$ws.Range("E19").Value = "  -6.18%  "
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("E21").Value = "  -7.62%  "
$ws.Range("D22").Value = "27.625.05"
$ws.Range("E22").Value = "  -4.44%  "
$ws.Range("D23").Value = "4.958"
$ws.Range("E23").Value = "  -6.95%  "
$ws.Range("E24").Value = "  -3.27%  "
$ws.Range("D25").Value = "2.067.96"
$ws.Range("E25").Value = "  -3.08%  "
$ws.Range("D26").Value = "2.050"
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").Value = "153.01"
$ws.Range("E27").Value = "  -2.89%  "
$ws.Range("D28").Value = "18.18"
$ws.Range("E28").Value = "  -4.46%  "
$ws.Range("D29").Value = "5.334"
$ws.Range("E29").Value = "  -4.78%  "
$ws.Range("D30").Value = "111.32"
$ws.Range("E30").Value = "  -5.55%  "
$ws.Range("D31").Value = "1.766"
$ws.Range("E31").Value = "  -3.92%  "
$ws.Range("D32").Value = "0.08906"
$ws.Range("E32").Value = "  -4.18%  "
$ws.Range("D33").Value = "0.7828"
$ws.Range("E33").Value = "  -9.62%  "
$ws.Range("D34").Value = "4.507"
$ws.Range("E34").Value = "  -11.45%  "
$ws.Range("D35").Value = "2.854"
$ws.Range("E35").Value = "  -5.29%  "
$ws.Range("D36").Value = "1.082"
$ws.Range("E36").Value = "  -12.64%  "
$ws.Range("D37").Value = "0.9996"
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("D38").Value = "0.05441"
$ws.Range("E38").Value = "  -4.22%  "
$ws.Range("E39").Value = "  -4.72%  "
$ws.Range("D40").Value = "2.984"
$ws.Range("E40").Value = "  -2.69%  "
$ws.Range("D41").Value = "0.01926"
$ws.Range("E41").Value = "  -5.42%  "
$ws.Range("E42").Value = "  -7.95%  "
$ws.Range("D43").Value = "6.779"
$ws.Range("E43").Value = "  -9.03%  "
$ws.Range("D44").Value = "0.1634"
$ws.Range("E44").Value = "  -6.73%  "
$ws.Range("D45").Value = "8.266"
$ws.Range("E45").Value = "  -11.14%  "
$ws.Range("D46").Value = "0.06622"
$ws.Range("E46").Value = "  -4.44%  "
$ws.Range("D47").Value = "106.53"
$ws.Range("E47").Value = "  -3.57%  "
$ws.Range("D48").Value = "0.4709"
$ws.Range("E48").Value = "  -8.53%  "
$ws.Range("D49").Value = "10.24"
$ws.Range("E49").Value = "  -7.59%  "
$ws.Range("D50").Value = "0.9992"
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("D51").Value = "1.646"
$ws.Range("E51").Value = "  -6.42%  "
